$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: becomes old row 9's data (Cylinder=4 group) ---
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 83.33333333333333
$ws.Range("E8").Value = 18.50225211517056
$ws.Range("F8").Value = 2.886666666666667
$ws.Range("G8").Value = 0.4911551010967242

# --- Row 9: becomes old row 8's data (Cylinder=6, Engine=0 group) ---
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 110
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 2.7475
$ws.Range("G9").Value = 0.1803122292025695

# --- Merge A9:A10 together to visually span the Cylinder=6 rows ---
# (do this before finalising formatting -- Merge() re-stamps the format of
#  every cell in the range using the anchor's style, so any formatting we
#  want to stick needs to be (re)applied afterwards)
$ws.Range("A9:A10").Merge()

# --- Re-establish A9's base look (same as A8, a Cylinder data cell) and then
#     layer the new top-vertical alignment on top of it ---
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").VerticalAlignment = -4160

# --- Row 10 stays Cylinder=6, Engine=1 group, but column A is now merged with A9,
#     so A10 goes back to looking like a plain/blank bordered cell (matches column H) ---
$ws.Range("H10").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").ClearContents()

# Row 11 and everything else remain untouched.
